$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The original table was a narrow A1:D5 block (Year in column A, three
# data columns, with a couple of blank spacer rows). The new table is a
# wider A1:H4 grid: Year stays in column A (keeping its existing bold +
# bordered style), and six new metric columns (with two "section title"
# columns B/F that only have a header, no per-row values) replace the
# old B/C/D columns.
# ---------------------------------------------------------------------

# Clear everything except column A's existing formatting, which already
# has the bold/bordered style (xf index 1) applied on A1/A2 that we want
# to reuse/extend down to A3:A4.
$ws.Range("B1:H5").Clear()
$ws.Range("A3:A5").Clear()

# Extend the existing "Year" header/style down through A3:A4 by copying
# formats only from A1 (format-only copy re-uses the same underlying
# style instead of creating a near-duplicate one).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A3:A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

function Set-TextValue($cellAddress, $text) {
    # Use a scratch cell far outside the target range to stage the value
    # as genuine text (via a "@" text number format) and then paste just
    # the *value* into the destination cell, leaving the destination's
    # own formatting completely untouched. This avoids Excel silently
    # re-interpreting numeric-looking strings (e.g. "2019", "33,104,461")
    # as numbers, and avoids leaving stray custom number formats behind
    # on the destination cell's style.
    $scratch = $ws.Range("Z100")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy() | Out-Null
    $ws.Range($cellAddress).PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $excel.CutCopyMode = 0
    $scratch.Clear()
}

# --- Header row (row 1): Year / category headers -------------------------
Set-TextValue "A1" "Year"
Set-TextValue "B1" "Electricity consumption"
Set-TextValue "C1" "Total electricity consumption"
Set-TextValue "D1" "of which green electricity"
Set-TextValue "E1" "Individual electricity consumption (kWh per FTE)"
Set-TextValue "F1" "District heating consumption"
Set-TextValue "G1" "Total district heating consumption"
Set-TextValue "H1" "Individual heating consumption (kWh per FTE)"

# --- Row 2: 2019 data ------------------------------------------------------
Set-TextValue "A2" "2019"
Set-TextValue "C2" "33,104,461"
Set-TextValue "D2" "32,782,553"
Set-TextValue "E2" "7,116"
Set-TextValue "G2" "18,964,126"
Set-TextValue "H2" "4,077"

# --- Row 3: 2018 data ------------------------------------------------------
Set-TextValue "A3" "2018"
Set-TextValue "C3" "33,035,150"
Set-TextValue "D3" "33,005,705"
Set-TextValue "E3" "7,000"
Set-TextValue "G3" "18,124,104"
Set-TextValue "H3" "3,841"

# --- Row 4: 2017 data ------------------------------------------------------
Set-TextValue "A4" "2017"
Set-TextValue "C4" "32,208,132"
Set-TextValue "D4" "32,036,926"
Set-TextValue "E4" "6,734"
Set-TextValue "G4" "18,563,309"
Set-TextValue "H4" "3,881"
